$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells keep their exact text representation
# (Excel would otherwise auto-convert numeric-looking strings to numbers)

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '95.871.89'

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.666.39'
$ws.Range('E3').Value = '  +10.38%  '

$ws.Range('E4').Value = '  +0.00%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '241.85'
$ws.Range('E5').Value = '  +4.94%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '642.95'
$ws.Range('E6').Value = '  +4.94%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.47'
$ws.Range('E7').Value = '  +5.52%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.401'

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.999'
$ws.Range('E9').Value = '  -0.15%  '

$ws.Range('E10').Value = '  +5.53%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '3.663.66'
$ws.Range('E11').Value = '  +10.31%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '43.63'
$ws.Range('E12').Value = '  +2.84%  '

$ws.Range('E13').Value = '  +3.86%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.38'
$ws.Range('E14').Value = '  +4.27%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.358.21'
$ws.Range('E15').Value = '  +10.39%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '95.766.49'
$ws.Range('E16').Value = '  +4.50%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000256'
$ws.Range('E17').Value = '  +5.46%  '

$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.675.40'
$ws.Range('E18').Value = '  +10.23%  '

$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.41'
$ws.Range('E19').Value = '  +24.31%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '8.05'
$ws.Range('E20').Value = '  +0.28%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '18.79'
$ws.Range('E21').Value = '  +8.91%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '519.47'
$ws.Range('E22').Value = '  +5.83%  '

$ws.Range('B23').Value = 'Stellar'
$ws.Range('C23').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.484'
$ws.Range('E23').Value = '  +10.76%  '

$ws.Range('B24').Value = 'SuiNetwork'
$ws.Range('C24').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.44'
$ws.Range('E24').Value = '  +0.97%  '

$ws.Range('E25').Value = '  +9.62%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.81'
$ws.Range('E26').Value = '  +5.22%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '97.53'
$ws.Range('E27').Value = '  +5.87%  '

$ws.Range('E28').Value = '  +6.21%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.17'
$ws.Range('E29').Value = '  +21.89%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '11.66'

$ws.Range('E31').Value = '  +2.94%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.999'
$ws.Range('E32').Value = '  -0.11%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '32.80'
$ws.Range('E33').Value = '  +16.23%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.180'
$ws.Range('E34').Value = '  +4.76%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.997'
$ws.Range('E35').Value = '  -1.83%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.577'
$ws.Range('E36').Value = '  +9.93%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '561.63'
$ws.Range('E37').Value = '  -0.19%  '

$ws.Range('E38').Value = '  +9.56%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '7.86'
$ws.Range('E39').Value = '  +6.57%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.961'
$ws.Range('E40').Value = '  +11.10%  '

$ws.Range('E41').Value = '  +2.75%  '

$ws.Range('E42').Value = '  -0.09%  '

$ws.Range('E43').Value = '  +7.91%  '

$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0432'
$ws.Range('E44').Value = '  +4.95%  '

$ws.Range('B45').Value = 'ImmutableX'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.74'
$ws.Range('E45').Value = '  +4.20%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '23.72'
$ws.Range('E46').Value = '  +0.14%  '

$ws.Range('B47').Value = 'Stacks'
$ws.Range('C47').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.22'
$ws.Range('E47').Value = '  +5.65%  '

$ws.Range('B48').Value = 'OKB'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '54.49'
$ws.Range('E48').Value = '  +5.18%  '

$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '32.65'
$ws.Range('E49').Value = '  +45.58%  '

$ws.Range('B50').Value = 'Cosmos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.31'
$ws.Range('E50').Value = '  +4.41%  '

$ws.Range('E51').Value = '  -2.47%  '
